$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.516.11'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.742.34'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4430'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3522'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07404'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.076'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.902'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.073'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.28%  '
$ws.Range("D16").Value = '1.740.21'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06380'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.723'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("D23").Value = '27.556.77'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.099'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("D28").Value = '1.940.87'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '124.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.035'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.56%  '
$ws.Range("E31").Value = '  -5.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09090'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.653'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.369'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.40%  '
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("E36").Value = '  -5.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06024'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2059'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6236'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.889'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.184'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.373'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.706'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.80%  '
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5784'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("E50").Value = '  -4.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.43%  '
